{"js": "// Update the two-digit x two-digit multiplication prompts in the table.\n// Each \"AxB=\" cell text is unique in the document, so a direct\n// search-and-replace (matching the whole cell text, case/whitespace\n// sensitive) is unambiguous and safe for every occurrence.\nconst replacements = [\n  [\"76\u00d791=\", \"51\u00d734=\"],\n  [\"37\u00d796=\", \"52\u00d753=\"],\n  [\"65\u00d789=\", \"75\u00d731=\"],\n  [\"92\u00d782=\", \"54\u00d778=\"],\n  [\"84\u00d744=\", \"38\u00d795=\"],\n  [\"44\u00d799=\", \"75\u00d759=\"],\n  [\"86\u00d789=\", \"98\u00d762=\"],\n  [\"46\u00d719=\", \"91\u00d744=\"],\n  [\"57\u00d777=\", \"20\u00d766=\"],\n  [\"54\u00d722=\", \"95\u00d732=\"],\n  [\"16\u00d757=\", \"18\u00d723=\"],\n  [\"23\u00d775=\", \"46\u00d755=\"],\n  [\"74\u00d718=\", \"47\u00d794=\"],\n  [\"32\u00d799=\", \"84\u00d725=\"],\n  [\"45\u00d735=\", \"36\u00d750=\"],\n  [\"36\u00d765=\", \"37\u00d715=\"],\n  [\"58\u00d775=\", \"58\u00d780=\"],\n  [\"59\u00d720=\", \"39\u00d713=\"],\n  [\"28\u00d787=\", \"67\u00d738=\"],\n  [\"99\u00d740=\", \"15\u00d793=\"],\n  [\"48\u00d711=\", \"43\u00d755=\"],\n  [\"28\u00d793=\", \"35\u00d713=\"],\n  [\"85\u00d761=\", \"72\u00d721=\"],\n  [\"80\u00d785=\", \"14\u00d762=\"],\n  [\"16\u00d760=\", \"78\u00d750=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the two-digit x two-digit multiplication prompts in the table.\n# Each \"AxB=\" cell text is unique in the document, so Find/Replace with\n# MatchCase/whole-string matching is unambiguous for every occurrence.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"76\u00d791=\", \"51\u00d734=\"),\n    @(\"37\u00d796=\", \"52\u00d753=\"),\n    @(\"65\u00d789=\", \"75\u00d731=\"),\n    @(\"92\u00d782=\", \"54\u00d778=\"),\n    @(\"84\u00d744=\", \"38\u00d795=\"),\n    @(\"44\u00d799=\", \"75\u00d759=\"),\n    @(\"86\u00d789=\", \"98\u00d762=\"),\n    @(\"46\u00d719=\", \"91\u00d744=\"),\n    @(\"57\u00d777=\", \"20\u00d766=\"),\n    @(\"54\u00d722=\", \"95\u00d732=\"),\n    @(\"16\u00d757=\", \"18\u00d723=\"),\n    @(\"23\u00d775=\", \"46\u00d755=\"),\n    @(\"74\u00d718=\", \"47\u00d794=\"),\n    @(\"32\u00d799=\", \"84\u00d725=\"),\n    @(\"45\u00d735=\", \"36\u00d750=\"),\n    @(\"36\u00d765=\", \"37\u00d715=\"),\n    @(\"58\u00d775=\", \"58\u00d780=\"),\n    @(\"59\u00d720=\", \"39\u00d713=\"),\n    @(\"28\u00d787=\", \"67\u00d738=\"),\n    @(\"99\u00d740=\", \"15\u00d793=\"),\n    @(\"48\u00d711=\", \"43\u00d755=\"),\n    @(\"28\u00d793=\", \"35\u00d713=\"),\n    @(\"85\u00d761=\", \"72\u00d721=\"),\n    @(\"80\u00d785=\", \"14\u00d762=\"),\n    @(\"16\u00d760=\", \"78\u00d750=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
